$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-30 Wednesday" "2024-10-31 Thursday"
Replace-Text "63×16=" "52×19="
Replace-Text "43×30=" "33×91="
Replace-Text "51×26=" "86×62="
Replace-Text "87×36=" "57×59="
Replace-Text "50×33=" "50×85="
Replace-Text "80×15=" "17×30="
Replace-Text "46×97=" "12×88="
Replace-Text "61×96=" "93×13="
Replace-Text "34×61=" "15×39="
Replace-Text "79×47=" "20×26="
Replace-Text "37×40=" "31×98="
Replace-Text "17×16=" "54×56="
Replace-Text "86×97=" "93×35="
Replace-Text "90×55=" "80×11="
Replace-Text "93×14=" "84×34="
Replace-Text "33×55=" "27×23="
Replace-Text "58×11=" "47×98="
Replace-Text "35×99=" "15×95="
Replace-Text "96×22=" "62×53="
Replace-Text "26×21=" "63×59="
Replace-Text "22×61=" "37×59="
Replace-Text "23×29=" "26×77="
Replace-Text "54×31=" "41×40="
Replace-Text "13×89=" "25×73="
Replace-Text "90×44=" "55×77="
